$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The filing on row 27 (index 25, dated 2018-05-15) is the last one that
# actually occurred - there were "no more filings after 2019-05-10". Rows 28
# through 52 no longer have a sequential filing index, so clear column A for
# those rows. Row 27 becomes the new "last" entry and its border loses the
# bottom edge (the box that used to close at row 28+ now closes at row 27,
# open at the bottom).

$ws.Range("A28:A52").ClearContents()

# Row 27: keep the bold/centered formatting but drop the bottom border
# (top+left+right thin border only).
$c27 = $ws.Cells.Item(27, 1)
$c27.Borders.LineStyle = 1             # xlContinuous, all four edges thin
$c27.Borders.Weight = 2                # xlThin
$c27.Borders.Item(9).LineStyle = -4142 # xlEdgeBottom -> none

# Rows 28-52: keep bold/centered formatting but remove all borders entirely.
# Apply cell-by-cell, and set the whole Borders collection in one shot so
# every cell individually resolves to the same "no border" style (instead of
# minting a new style per edge-assignment / per range-boundary quirk).
for ($r = 28; $r -le 52; $r++) {
    $ws.Cells.Item($r, 1).Borders.LineStyle = -4142
}

# Restore the sheet view: scroll back to the top and move the selection.
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("G23").Select()
